# Replace the legal-basis enumeration paragraphs (items "2." through "6.")
# with a new item "2." (ΠΔ 18/2018) inserted before the rest, renumbering the
# remaining items 3..7 accordingly (old 2->3, 3->4, 4->5, 5->6, 6->7).

$d = $word.ActiveDocument

# Locate the start of the block: the paragraph that begins with old item "2."
# (unique text fragment "Φ.353.1/324/105657").
$startRng = $d.Content
$startRng.Find.Execute('Φ.353.1/324/105657', $true, $false, $false, $false, $false, $true, 1, $false, '', 0) | Out-Null
$startRng.Expand(4) | Out-Null

# Locate the end of the block: the paragraph that ends old item "6."
# (unique text fragment "με τα συνημμένα δικαιολογητικά").
$endRng = $d.Content
$endRng.Find.Execute('με τα συνημμένα δικαιολογητικά', $true, $false, $false, $false, $false, $true, 1, $false, '', 0) | Out-Null
$endRng.Expand(4) | Out-Null

$blockRng = $d.Range($startRng.Start, $endRng.End)

$newXml = '<w:p><w:pPr><w:spacing w:after="60"/><w:ind w:left="284" w:hanging="284"/><w:jc w:val="both"/></w:pPr><w:r><w:t>2. Το ΠΔ 18/2018 (ΦΕΚ 31/τ.Α’/23-03-2018) «Οργανισμός Υπουργείου Παιδείας, Έρευνας και Θρησκευμάτων»</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="60"/><w:ind w:left="283" w:hanging="283"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>Τη με αριθ. Φ.353.1/324/105657/Δ1 (ΦΕΚ 1340/16-10-2002 τ Β’) απόφαση του Υπουργού Εθνικής Παιδείας και Θρησκευμάτων με θέμα «</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i/></w:rPr><w:t>Καθορισμός καθηκόντων και αρμοδιοτήτων των προϊσταμένων των περιφερειακών υπηρεσιών πρωτοβάθμιας και δευτεροβάθμιας εκπαίδευσης, των διευθυντών και υποδιευθυντών των σχολικών μονάδων και ΣΕΚ και των συλλόγων των διδασκόντων</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>», όπως τροποποιήθηκε, συμπληρώθηκε και ισχύει.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="60"/><w:ind w:left="283" w:hanging="283"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Τη με αριθ. Φ.350.2/1/32958/E3/27-02-2018 (ΑΔΑ:6Π414653ΠΣ-7ΕΝ) Υπουργική Απόφαση με θέμα «</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>Τοποθέτηση Περιφερειακών Διευθυντών Εκπαίδευσης</w:t></w:r><w:r><w:t>».</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="60"/><w:ind w:left="283" w:hanging="283"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>5</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>Τη με αριθ. ${PLACEMENT_NUM} και ημερομηνία ${PLAC_DATE} απόφαση τοποθέτησης με θέμα: «</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:i/></w:rPr><w:t>${PLAC_SUBJ}</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>».</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="60"/><w:ind w:left="283" w:hanging="283"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>6</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>. Την ανάγκη μετάβασης εκτός έδρας για εκτέλεση υπηρεσίας.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="60"/><w:ind w:left="283" w:hanging="283"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>7</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:szCs w:val="24"/></w:rPr><w:t>${APPLIC}</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> με τα συνημμένα δικαιολογητικά.</w:t></w:r></w:p>'

$blockRng.InsertXML($newXml)
